$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 495.33334
$ws.Range("I28").Value = 495.33334
$ws.Range("K28").Value = 495.33334
$ws.Range("M28").Value = -10.33334000000002

$ws.Range("H76").Value = 19239840
$ws.Range("I76").Value = 10284.363
$ws.Range("K76").Value = 10284.363
$ws.Range("M76").Value = -9969.362999999999

$ws.Range("H79").Value = 19239840
$ws.Range("I79").Value = 10284.363
$ws.Range("K79").Value = 10284.363
$ws.Range("M79").Value = -9192.362999999999

$ws.Range("H92").Value = 2872.875
$ws.Range("J92").Value = 5249
$ws.Range("L92").Value = 5249
$ws.Range("N92").Value = -7745

$ws.Range("I107").Value = 13890126
$ws.Range("J107").Value = 125000250
$ws.Range("K107").Value = 13890126
$ws.Range("L107").Value = 125000250
$ws.Range("M107").Value = -13888206
$ws.Range("N107").Value = -125004090

$ws.Range("H118").Value = 1449.1666
$ws.Range("I118").Value = 692.44446
$ws.Range("J118").Value = 3719.3333
$ws.Range("K118").Value = 2077.33338
$ws.Range("L118").Value = 11157.9999
$ws.Range("M118").Value = -420.33338
$ws.Range("N118").Value = -14471.9999

$ws.Range("H133").Value = 100770
$ws.Range("J133").Value = 100770
$ws.Range("L133").Value = 100770
$ws.Range("N133").Value = -110890

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 4368.6665
$ws.Range("J14").Value = 4368.6665
$ws.Range("L14").Value = 4368.6665
$ws.Range("N14").Value = -4718.6665

$ws.Range("H32").Value = 2845364.5
$ws.Range("I32").Value = 3294014.2
$ws.Range("K32").Value = 3294014.2
$ws.Range("M32").Value = -3293727.2

$ws.Range("H45").Value = 4281.72
$ws.Range("I45").Value = 1120.3846
$ws.Range("K45").Value = 1120.3846
$ws.Range("M45").Value = -743.3846000000001

$ws.Range("H97").Value = 2782426.2
$ws.Range("I97").Value = 684.6316
$ws.Range("K97").Value = 684.6316
$ws.Range("M97").Value = -188.6316

$ws.Range("H98").Value = 54445
$ws.Range("J98").Value = 54445
$ws.Range("L98").Value = 54445
$ws.Range("N98").Value = -60435

$ws.Range("H102").Value = 4200.4
$ws.Range("I102").Value = 4151.75
$ws.Range("J102").Value = 4395
$ws.Range("K102").Value = 4151.75
$ws.Range("L102").Value = 4395
$ws.Range("M102").Value = -2529.75
$ws.Range("N102").Value = -7639

$ws.Range("I110").Value = 2316.7778
$ws.Range("K110").Value = 2316.7778
$ws.Range("M110").Value = -271.7777999999998

$ws.Range("H122").Value = 3323.5789
$ws.Range("I122").Value = 2687.9355
$ws.Range("K122").Value = 8063.806500000001
$ws.Range("M122").Value = -5613.806500000001

$ws.Range("H132").Value = 6817.654
$ws.Range("I132").Value = 3770.9333
$ws.Range("J132").Value = 10972.272
$ws.Range("K132").Value = 11312.7999
$ws.Range("L132").Value = 32916.81600000001
$ws.Range("M132").Value = -8782.7999
$ws.Range("N132").Value = -37976.81600000001

$ws.Range("H139").Value = 82000
$ws.Range("J139").Value = 82000
$ws.Range("L139").Value = 82000
$ws.Range("N139").Value = -92280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -887

$ws.Range("H20").Value = 6667876.5
$ws.Range("J20").Value = 947.7778
$ws.Range("L20").Value = 947.7778
$ws.Range("N20").Value = -1441.7778

$ws.Range("H86").Value = 34787.97
$ws.Range("I86").Value = 57157.11
$ws.Range("K86").Value = 57157.11
$ws.Range("M86").Value = -56034.11

$ws.Range("H89").Value = 34787.97
$ws.Range("I89").Value = 57157.11
$ws.Range("K89").Value = 285785.55
$ws.Range("M89").Value = -280169.55

$ws.Range("H105").Value = 2651.1052
$ws.Range("I105").Value = 1800
$ws.Range("K105").Value = 1800
$ws.Range("M105").Value = -53

$ws.Range("H134").Value = 10062.611
$ws.Range("I134").Value = 4607.25
$ws.Range("J134").Value = 11621.286
$ws.Range("K134").Value = 13821.75
$ws.Range("L134").Value = 34863.858
$ws.Range("M134").Value = -11286.75
$ws.Range("N134").Value = -39933.858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 40000
$ws.Range("J100").Value = 40000
$ws.Range("L100").Value = 40000
$ws.Range("N100").Value = -42164

$ws.Range("H122").Value = 3793.9333
$ws.Range("I122").Value = 3637.2727
$ws.Range("K122").Value = 10911.8181
$ws.Range("M122").Value = -8461.8181

$ws.Range("H134").Value = 9039.289000000001
$ws.Range("I134").Value = 9697.883
$ws.Range("J134").Value = 8506.143
$ws.Range("K134").Value = 29093.649
$ws.Range("L134").Value = 25518.429
$ws.Range("M134").Value = -26558.649
$ws.Range("N134").Value = -30588.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1999.5
$ws.Range("I3").Value = 1999.5
$ws.Range("K3").Value = 5998.5
$ws.Range("M3").Value = -5886.5

$ws.Range("H68").Value = 7238.5713
$ws.Range("I68").Value = 1649.6666
$ws.Range("K68").Value = 4948.9998
$ws.Range("M68").Value = -4137.9998

$ws.Range("H71").Value = 7238.5713
$ws.Range("I71").Value = 1649.6666
$ws.Range("K71").Value = 14846.9994
$ws.Range("M71").Value = -10790.9994

$ws.Range("H107").Value = 813.5
$ws.Range("I107").Value = 633
$ws.Range("K107").Value = 1899
$ws.Range("M107").Value = 21

$ws.Range("H122").Value = 5660265
$ws.Range("I122").Value = 9430109
$ws.Range("K122").Value = 84870981
$ws.Range("M122").Value = -84868531

$ws.Range("H131").Value = 1851.2222
$ws.Range("I131").Value = 1002.8461
$ws.Range("J131").Value = 4057
$ws.Range("K131").Value = 3008.5383
$ws.Range("L131").Value = 12171
$ws.Range("M131").Value = 2031.4617
$ws.Range("N131").Value = -22251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1250182
$ws.Range("I2").Value = 102.4
$ws.Range("J2").Value = 3333648.2
$ws.Range("K2").Value = 102.4
$ws.Range("L2").Value = 3333648.2
$ws.Range("M2").Value = 10.59999999999999
$ws.Range("N2").Value = -3333874.2

$ws.Range("H52").Value = 48744
$ws.Range("J52").Value = 48744
$ws.Range("L52").Value = 48744
$ws.Range("N52").Value = -49262

$ws.Range("H57").Value = 49996.445
$ws.Range("J57").Value = 49996.445
$ws.Range("L57").Value = 49996.445
$ws.Range("N57").Value = -51636.445

$ws.Range("H58").Value = 54397.285
$ws.Range("J58").Value = 54397.285
$ws.Range("L58").Value = 54397.285
$ws.Range("N58").Value = -54951.285

$ws.Range("H126").Value = 62502476
$ws.Range("I126").Value = 166669140
$ws.Range("J126").Value = 2479.4
$ws.Range("K126").Value = 500007420
$ws.Range("L126").Value = 7438.200000000001
$ws.Range("M126").Value = -500004950
$ws.Range("N126").Value = -12378.2

$ws.Range("H132").Value = 8103.8184
$ws.Range("I132").Value = 2028.4
$ws.Range("K132").Value = 6085.200000000001
$ws.Range("M132").Value = -3555.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5294.84
$ws.Range("I22").Value = 2165.8823
$ws.Range("J22").Value = 11943.875
$ws.Range("K22").Value = 2165.8823
$ws.Range("L22").Value = 11943.875
$ws.Range("M22").Value = -1870.8823
$ws.Range("N22").Value = -12533.875

$ws.Range("H27").Value = 5294.84
$ws.Range("I27").Value = 2165.8823
$ws.Range("J27").Value = 11943.875
$ws.Range("K27").Value = 2165.8823
$ws.Range("L27").Value = 11943.875
$ws.Range("M27").Value = -2058.8823
$ws.Range("N27").Value = -12157.875

$ws.Range("H46").Value = 1511.5769
$ws.Range("I46").Value = 810.875
$ws.Range("K46").Value = 810.875
$ws.Range("M46").Value = -622.875

$ws.Range("H55").Value = 332.8889
$ws.Range("J55").Value = 441.92307
$ws.Range("L55").Value = 441.92307
$ws.Range("N55").Value = -787.9230700000001

$ws.Range("H61").Value = 3956.8484
$ws.Range("I61").Value = 2884.1667
$ws.Range("J61").Value = 5244.067
$ws.Range("K61").Value = 2884.1667
$ws.Range("L61").Value = 5244.067
$ws.Range("M61").Value = -2682.1667
$ws.Range("N61").Value = -5648.067

$ws.Range("H100").Value = 2670.6428
$ws.Range("I100").Value = 2378.4443
$ws.Range("K100").Value = 2378.4443
$ws.Range("M100").Value = -1837.4443

$ws.Range("H113").Value = 3956.8484
$ws.Range("I113").Value = 2884.1667
$ws.Range("J113").Value = 5244.067
$ws.Range("K113").Value = 2884.1667
$ws.Range("L113").Value = 5244.067
$ws.Range("M113").Value = -714.1667000000002
$ws.Range("N113").Value = -9584.066999999999

$ws.Range("H132").Value = 8000.1724
$ws.Range("I132").Value = 4066.6667
$ws.Range("K132").Value = 12200.0001
$ws.Range("M132").Value = -9670.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

$ws.Range("H122").Value = 14825998
$ws.Range("I122").Value = 20161614
$ws.Range("K122").Value = 60484842
$ws.Range("M122").Value = -60482392

$ws.Range("H125").Value = 51644
$ws.Range("J125").Value = 51644
$ws.Range("L125").Value = 51644
$ws.Range("N125").Value = -61484

$ws.Range("H136").Value = 38031.137
$ws.Range("I136").Value = 1975.409
$ws.Range("K136").Value = 5926.227000000001
$ws.Range("M136").Value = -3376.227000000001
